# ------------------------------------------------------------------
# This script updates the OpenSea "orders" sheet so that it contains
# four data rows instead of two:
#   - row 2 and row 3 are refreshed with newer scrape values
#   - two brand new rows (4 and 5) are appended
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 4-5 do not exist yet, so Excel will not know they should
# use the bordered/centered "index" style (s="1") that column A uses
# on every other data row. Copy that formatting from A3 first.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Row 2: update only the columns that actually changed ----
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2021-07-23T02:26:51.337998"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2021-07-23T04:26:46"
$ws.Range("E2").Value = 1627014406
$ws.Range("F2").Value = 1627007107
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0x3724a6bad890de8998aed8c6fcc61e5a7dfcd59d953edebee6a53c33495b0e5c"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "133000000000000000.0000000000"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "1330000000000000"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "0x23b872dd00000000000000000000000000000000000000000000000000000000000000000000000000000000000000008249b4417d752cda7022554b111a2f7c3323da060000000000000000000000000000000000000000000000000000000000001566"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "133000000000000000"
$ws.Range("AD2").NumberFormat = "@"
$ws.Range("AD2").Value = "20639483450222611293504570282524924318059802967165126658734627575340007104385"
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = "0x59afbdadffc381136184b8dba5373a1cbc830ba0eb9c5c3cee9daba2a5cb5c88"
$ws.Range("AG2").NumberFormat = "@"
$ws.Range("AG2").Value = "0x446cad6dde2c0e4626cb6ca5e167ad04463e48778214476bc7436d5c597444d6"
$ws.Range("AL2").NumberFormat = "@"
$ws.Range("AL2").Value = "0xcb710e83254751118ab63c7bfd620bbae52ce9989fa8d2b9b003a447d537686f"
$ws.Range("AP2").NumberFormat = "@"
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").NumberFormat = "@"
$ws.Range("AQ2").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/7.png"
$ws.Range("AR2").NumberFormat = "@"
$ws.Range("AR2").Value = "0x8249b4417d752cda7022554b111a2f7c3323da06"
$ws.Range("BL2").NumberFormat = "@"
$ws.Range("BL2").Value = "2063.719999999999800000"

# ---- Row 3: update only the columns that actually changed ----
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2021-07-23T00:54:15.905532"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2021-07-23T23:21:19"
$ws.Range("E3").Value = 1627082479
$ws.Range("F3").Value = 1627001553
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0x0d0c7d723f058493c99c90148cad7a5002cd10acf8518b4e5dbdbdc25f998cc7"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "130000000000000000.0000000000"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "1300000000000000"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = "0x23b872dd0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000792a4ba75d04ca9ca369328736844e64928080920000000000000000000000000000000000000000000000000000000000001566"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "130000000000000000"
$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "18847090867684994452233122722601024441187212529845096456333665889826023711380"
$ws.Range("AF3").NumberFormat = "@"
$ws.Range("AF3").Value = "0x5d9d0c3ee6dd128aba0788610910d1a26cb6319ea5caa50b5dafa32b9bd40f8f"
$ws.Range("AG3").NumberFormat = "@"
$ws.Range("AG3").Value = "0x731176f790b23aecbdfe7e5d110a357782662a295982ae59b65928aa15b91f66"
$ws.Range("AL3").NumberFormat = "@"
$ws.Range("AL3").Value = "0x924e8e91a363527a0275e8dcfdfa2c15f9b9808f5fa370ba626ff8f35c37c17a"
$ws.Range("AP3").NumberFormat = "@"
$ws.Range("AP3").Value = "Butters666"
$ws.Range("AR3").NumberFormat = "@"
$ws.Range("AR3").Value = "0x792a4ba75d04ca9ca369328736844e6492808092"
$ws.Range("BL3").NumberFormat = "@"
$ws.Range("BL3").Value = "2063.719999999999800000"

# ---- Row 4: brand new row, set every column ----
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2021-07-22T21:41:13.013857"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2021-07-23T20:49:31"
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = 1627073371
$ws.Range("F4").Value = 1626989965
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0x892b3e5cd161c16ce9afb0f376d2695e5c24b3ce0ad0e54f0a2b24b24c6c6f14"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0x7be8076f4ea4a4ad08075c2508e481d6c946d12b"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "125000000000000000.0000000000"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "1250000000000000"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "0.01"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "0"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "550"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "0"
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "0"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "0"
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = "0x3fe1a4c1481c8351e91b64d5c398b159de07cbc5"
$ws.Range("U4").Value = 0
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value = "0x23b872dd000000000000000000000000000000000000000000000000000000000000000000000000000000000000000020cc3f4cb2df386f0b82713f0cc5464bfa6c05740000000000000000000000000000000000000000000000000000000000001566"
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value = "0x00000000ffffffffffffffffffffffffffffffffffffffffffffffffffffffffffffffff00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000"
$ws.Range("X4").NumberFormat = "@"
$ws.Range("X4").Value = "0x0000000000000000000000000000000000000000"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "0x"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "0xc02aaa39b223fe8d0a0e5c4f27ead9083c756cc2"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "125000000000000000"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "0"
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = "1"
$ws.Range("AD4").NumberFormat = "@"
$ws.Range("AD4").Value = "94615897778944020430428835460705028788793008878645173732692119753802041938409"
$ws.Range("AE4").Value = 27
$ws.Range("AF4").NumberFormat = "@"
$ws.Range("AF4").Value = "0x31f7b8a330f3091332b28c9bf9c48610a3c27a25e3d00c499d3345d33f68786c"
$ws.Range("AG4").NumberFormat = "@"
$ws.Range("AG4").Value = "0x442e954f6537ba7765a19cc92ec8482a9cfa09116446440063b8f44f3b5a67fb"
$ws.Range("AH4").Value = $false
$ws.Range("AI4").Value = $false
$ws.Range("AJ4").Value = $false
$ws.Range("AK4").Value = $false
$ws.Range("AL4").NumberFormat = "@"
$ws.Range("AL4").Value = "0x7a83c7715aeb54e9611140ef1abc6a92e81fe4a428c2021791b02880bb8283f3"
$ws.Range("AM4").NumberFormat = "@"
$ws.Range("AM4").Value = "5478"
$ws.Range("AN4").NumberFormat = "@"
$ws.Range("AN4").Value = "0x3fe1a4c1481c8351e91b64d5c398b159de07cbc5"
$ws.Range("AO4").NumberFormat = "@"
$ws.Range("AO4").Value = "ERC721"
$ws.Range("AP4").NumberFormat = "@"
$ws.Range("AP4").Value = "Sanders333"
$ws.Range("AQ4").NumberFormat = "@"
$ws.Range("AQ4").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/21.png"
$ws.Range("AR4").NumberFormat = "@"
$ws.Range("AR4").Value = "0x20cc3f4cb2df386f0b82713f0cc5464bfa6c0574"
$ws.Range("AS4").NumberFormat = "@"
$ws.Range("AS4").Value = ""
$ws.Range("AT4").NumberFormat = "@"
$ws.Range("AT4").Value = ""
$ws.Range("AU4").NumberFormat = "@"
$ws.Range("AU4").Value = "NullAddress"
$ws.Range("AV4").NumberFormat = "@"
$ws.Range("AV4").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/1.png"
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value = "0x0000000000000000000000000000000000000000"
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value = ""
$ws.Range("AY4").NumberFormat = "@"
$ws.Range("AY4").Value = ""
$ws.Range("AZ4").NumberFormat = "@"
$ws.Range("AZ4").Value = "OS-Wallet"
$ws.Range("BA4").NumberFormat = "@"
$ws.Range("BA4").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/28.png"
$ws.Range("BB4").NumberFormat = "@"
$ws.Range("BB4").Value = "0x5b3256965e7c3cf26e11fcaf296dfc8807c01073"
$ws.Range("BC4").NumberFormat = "@"
$ws.Range("BC4").Value = "verified"
$ws.Range("BD4").NumberFormat = "@"
$ws.Range("BD4").Value = ""
$ws.Range("BE4").Value = 2
$ws.Range("BF4").NumberFormat = "@"
$ws.Range("BF4").Value = "WETH"
$ws.Range("BG4").NumberFormat = "@"
$ws.Range("BG4").Value = "0xc02aaa39b223fe8d0a0e5c4f27ead9083c756cc2"
$ws.Range("BH4").NumberFormat = "@"
$ws.Range("BH4").Value = "https://storage.opensea.io/files/accae6b6fb3888cbff27a013729c22dc.svg"
$ws.Range("BI4").NumberFormat = "@"
$ws.Range("BI4").Value = "Wrapped Ether"
$ws.Range("BJ4").Value = 18
$ws.Range("BK4").NumberFormat = "@"
$ws.Range("BK4").Value = "1.000000000000000"
$ws.Range("BL4").NumberFormat = "@"
$ws.Range("BL4").Value = "2063.719999999999800000"

# ---- Row 5: brand new row, set every column ----
$ws.Range("A5").Value = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2021-07-22T13:15:17.264875"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2021-07-23T10:19:42"
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = 1627035582
$ws.Range("F5").Value = 1626959612
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0xef510cdae2f4d94b153a6b71eca83bd81ce0c863b7841d13e10b54b0ddbe8de0"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0x7be8076f4ea4a4ad08075c2508e481d6c946d12b"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "123000000000000000.0000000000"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "1230000000000000"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "0.01"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "0"
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "550"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "0"
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "0"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "0"
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = "0x3fe1a4c1481c8351e91b64d5c398b159de07cbc5"
$ws.Range("U5").Value = 0
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value = "0x23b872dd000000000000000000000000000000000000000000000000000000000000000000000000000000000000000020cc3f4cb2df386f0b82713f0cc5464bfa6c05740000000000000000000000000000000000000000000000000000000000001566"
$ws.Range("W5").NumberFormat = "@"
$ws.Range("W5").Value = "0x00000000ffffffffffffffffffffffffffffffffffffffffffffffffffffffffffffffff00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000"
$ws.Range("X5").NumberFormat = "@"
$ws.Range("X5").Value = "0x0000000000000000000000000000000000000000"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "0x"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "0xc02aaa39b223fe8d0a0e5c4f27ead9083c756cc2"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "123000000000000000"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "0"
$ws.Range("AC5").NumberFormat = "@"
$ws.Range("AC5").Value = "1"
$ws.Range("AD5").NumberFormat = "@"
$ws.Range("AD5").Value = "40894784561982867332627611566104062480602654460596944182411762954639387733994"
$ws.Range("AE5").Value = 27
$ws.Range("AF5").NumberFormat = "@"
$ws.Range("AF5").Value = "0xec2934828c5ebb8354975fd34b5f3fd51f51e230dd6fe8e97aa78427e5857763"
$ws.Range("AG5").NumberFormat = "@"
$ws.Range("AG5").Value = "0x51cf8bf01f1c686bd8e5305214f69115b50dc0b373f12e9ebd514a29ec2eea23"
$ws.Range("AH5").Value = $false
$ws.Range("AI5").Value = $false
$ws.Range("AJ5").Value = $false
$ws.Range("AK5").Value = $false
$ws.Range("AL5").NumberFormat = "@"
$ws.Range("AL5").Value = "0xba53282646dc38c15244d552e74929520c75411ac688c3cd198fd49c5f8cea34"
$ws.Range("AM5").NumberFormat = "@"
$ws.Range("AM5").Value = "5478"
$ws.Range("AN5").NumberFormat = "@"
$ws.Range("AN5").Value = "0x3fe1a4c1481c8351e91b64d5c398b159de07cbc5"
$ws.Range("AO5").NumberFormat = "@"
$ws.Range("AO5").Value = "ERC721"
$ws.Range("AP5").NumberFormat = "@"
$ws.Range("AP5").Value = "Sanders333"
$ws.Range("AQ5").NumberFormat = "@"
$ws.Range("AQ5").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/21.png"
$ws.Range("AR5").NumberFormat = "@"
$ws.Range("AR5").Value = "0x20cc3f4cb2df386f0b82713f0cc5464bfa6c0574"
$ws.Range("AS5").NumberFormat = "@"
$ws.Range("AS5").Value = ""
$ws.Range("AT5").NumberFormat = "@"
$ws.Range("AT5").Value = ""
$ws.Range("AU5").NumberFormat = "@"
$ws.Range("AU5").Value = "NullAddress"
$ws.Range("AV5").NumberFormat = "@"
$ws.Range("AV5").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/1.png"
$ws.Range("AW5").NumberFormat = "@"
$ws.Range("AW5").Value = "0x0000000000000000000000000000000000000000"
$ws.Range("AX5").NumberFormat = "@"
$ws.Range("AX5").Value = ""
$ws.Range("AY5").NumberFormat = "@"
$ws.Range("AY5").Value = ""
$ws.Range("AZ5").NumberFormat = "@"
$ws.Range("AZ5").Value = "OS-Wallet"
$ws.Range("BA5").NumberFormat = "@"
$ws.Range("BA5").Value = "https://storage.googleapis.com/opensea-static/opensea-profile/28.png"
$ws.Range("BB5").NumberFormat = "@"
$ws.Range("BB5").Value = "0x5b3256965e7c3cf26e11fcaf296dfc8807c01073"
$ws.Range("BC5").NumberFormat = "@"
$ws.Range("BC5").Value = "verified"
$ws.Range("BD5").NumberFormat = "@"
$ws.Range("BD5").Value = ""
$ws.Range("BE5").Value = 2
$ws.Range("BF5").NumberFormat = "@"
$ws.Range("BF5").Value = "WETH"
$ws.Range("BG5").NumberFormat = "@"
$ws.Range("BG5").Value = "0xc02aaa39b223fe8d0a0e5c4f27ead9083c756cc2"
$ws.Range("BH5").NumberFormat = "@"
$ws.Range("BH5").Value = "https://storage.opensea.io/files/accae6b6fb3888cbff27a013729c22dc.svg"
$ws.Range("BI5").NumberFormat = "@"
$ws.Range("BI5").Value = "Wrapped Ether"
$ws.Range("BJ5").Value = 18
$ws.Range("BK5").NumberFormat = "@"
$ws.Range("BK5").Value = "1.000000000000000"
$ws.Range("BL5").NumberFormat = "@"
$ws.Range("BL5").Value = "2063.719999999999800000"

Write-Host "done"
